# [MOSIP-14336] Updating Masterdata utility
# Rework the loc_hierarchy_list worksheet:
#  - rename/reorder header columns (lang_code, hierarchy_level, hierarchy_level_name, is_active)
#  - drop the old "hierarchyLevel" header in column A
#  - replace sample data with per-language (eng/fra) hierarchy level rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row -------------------------------------------------------
# Copy the existing header style (index already used by B1/C1/D1) onto the
# new E1 header cell before we touch its value.
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)

$ws.Range("B1").Value = "lang_code"
$ws.Range("C1").Value = "hierarchy_level"
$ws.Range("D1").Value = "hierarchy_level_name"
$ws.Range("E1").Value = "is_active"

# The old column A header ("hierarchyLevel") is removed entirely.
$ws.Range("A1").Clear()

# ---- Data rows ----------------------------------------------------------
# id, lang_code, hierarchy_level, hierarchy_level_name, is_active
$data = @(
    @(0,  "eng", 0, "Country",     $true),
    @(1,  "fra", 0, "Pays",        $true),
    @(2,  "eng", 1, "Region",      $true),
    @(3,  "fra", 1, "Région",      $true),
    @(4,  "eng", 2, "Province",    $true),
    @(5,  "fra", 2, "Province",    $true),
    @(6,  "eng", 3, "City",        $true),
    @(7,  "fra", 3, "Ville",       $true),
    @(8,  "eng", 4, "Zone",        $true),
    @(9,  "fra", 4, "Zone",        $true),
    @(10, "eng", 5, "Postal Code", $true),
    @(11, "fra", 5, "code postal", $true)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rec = $data[$i]

    # Column A keeps the bold/centered/bordered header style, same as before.
    $ws.Range("B1").Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
}

Write-Output "loc_hierarchy_list updated"
